# edit.ps1 - Applies the "ISP Attribution" AnalyticPlan.docx revision
# described by the commit "started getClients, updated Analytic plan step 1"
#
# Summary of changes:
#  1. First paragraph ("ISP Attribution ") switches from manually-formatted
#     text to the built-in Heading1 style, plus a -450 twip left indent.
#  2. A bookmark named _Hlk112937288 is added around the "Dataset folder"
#     .. "isp_attr" paragraphs.
#  3. The lastRenderedPageBreak marker moves from the start of the
#     "Total: 8 NPI's..." run up to the start of the "Aug 2022" run.
#  4. Four separate runs that spell out "Split ID 3355 > Couldn't find NPI"
#     collapse into a single run.

$d = $word.ActiveDocument

function Set-RangeOpenXml {
    param(
        [__ComObject]$Range,
        [string]$BodyXml
    )
    $pkg = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
        '<w:body>' + $BodyXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    $Range.InsertXML($pkg)
}

function Get-ParagraphTextNoMark {
    param([__ComObject]$Paragraph)
    $t = $Paragraph.Range.Text
    # Paragraphs' Range.Text includes the trailing paragraph-mark character;
    # strip it so we can safely re-embed the text inside a fresh <w:t>.
    if ($t.Length -gt 0) {
        $lastCode = [int][char]$t.Substring($t.Length - 1)
        if ($lastCode -eq 13 -or $lastCode -eq 7) {
            $t = $t.Substring(0, $t.Length - 1)
        }
    }
    return $t
}

function Escape-Xml {
    param([string]$Text)
    return $Text.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
}

# ---------------------------------------------------------------------------
# 1) "ISP Attribution" paragraph -> Heading1 style with a -450 twip indent
# ---------------------------------------------------------------------------
$pTitle = $d.Paragraphs(1)
$pTitle.Style = "Heading 1"
$pTitle.LeftIndent = -22.5   # points; -450 twips / 20 = -22.5 pt

# ---------------------------------------------------------------------------
# 2) Bookmark around "Dataset folder" ... "isp_attr"
# ---------------------------------------------------------------------------
$pDatasetFolder = $null
$pIspAttr = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs($i).Range.Text
    if ($pDatasetFolder -eq $null -and $txt -like "*Dataset folder*") {
        $pDatasetFolder = $d.Paragraphs($i)
    }
    if ($txt -like "*isp_attr*") {
        $pIspAttr = $d.Paragraphs($i)
    }
}
if ($pDatasetFolder -ne $null -and $pIspAttr -ne $null) {
    $bmRange = $d.Range($pDatasetFolder.Range.Start, $pIspAttr.Range.End)
    $d.Bookmarks.Add("_Hlk112937288", $bmRange)
}

# ---------------------------------------------------------------------------
# 3) Move lastRenderedPageBreak from "Total: 8 NPI's..." to "Aug 2022"
# ---------------------------------------------------------------------------
$pAug = $null
$pTotal = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs($i).Range.Text
    if ($pAug -eq $null -and $txt -like "Aug 2022*") {
        $pAug = $d.Paragraphs($i)
    }
    if ($txt -like "Total: 8 NPI*") {
        $pTotal = $d.Paragraphs($i)
    }
}

if ($pAug -ne $null) {
    $augText = Get-ParagraphTextNoMark $pAug
    $augEsc = Escape-Xml $augText
    $augBody = '<w:p><w:pPr><w:pStyle w:val="Heading4"/><w:ind w:left="-810"/></w:pPr>' +
        '<w:r><w:lastRenderedPageBreak/><w:t xml:space="preserve">' + $augEsc + '</w:t></w:r></w:p>'
    Set-RangeOpenXml -Range $pAug.Range -BodyXml $augBody
}

if ($pTotal -ne $null) {
    $totalText = Get-ParagraphTextNoMark $pTotal
    $totalEsc = Escape-Xml $totalText
    $totalBody = '<w:p><w:pPr><w:ind w:left="-810" w:firstLine="0"/>' +
        '<w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/></w:rPr>' +
        '<w:t xml:space="preserve">' + $totalEsc + '</w:t></w:r></w:p>'
    Set-RangeOpenXml -Range $pTotal.Range -BodyXml $totalBody
}

# ---------------------------------------------------------------------------
# 4) Merge the four runs that make up "Split ID 3355 > Couldn't find NPI"
# ---------------------------------------------------------------------------
$pSplit = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $txt = $d.Paragraphs($i).Range.Text
    if ($txt -like "Split ID 3355*") {
        $pSplit = $d.Paragraphs($i)
        break
    }
}
if ($pSplit -ne $null) {
    $splitText = Get-ParagraphTextNoMark $pSplit
    $splitEsc = Escape-Xml $splitText
    $splitBody = '<w:p><w:pPr><w:pStyle w:val="ListParagraph"/>' +
        '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="7"/></w:numPr>' +
        '<w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/></w:rPr></w:pPr>' +
        '<w:r><w:rPr><w:rFonts w:ascii="Arial Narrow" w:hAnsi="Arial Narrow"/></w:rPr>' +
        '<w:t xml:space="preserve">' + $splitEsc + '</w:t></w:r></w:p>'
    Set-RangeOpenXml -Range $pSplit.Range -BodyXml $splitBody
}

Write-Host "Edits applied."
